# This script applies a reshuffle of the data rows (rows 2-12) in the
# "Artfynd" worksheet. Each destination row ends up with the full set of
# values (all columns) that a particular source row held before the
# change - i.e. the data rows get permuted / reordered while the header
# row (row 1) is left untouched.
#
# Mapping of destination row -> source row (both referring to the
# ORIGINAL, pre-edit layout):
#   2  <- 7
#   3  <- 11
#   4  <- 5
#   5  <- 4
#   6  <- 3
#   7  <- 8
#   8  <- 6
#   9  <- 10
#   10 <- 2
#   11 <- 12
#   12 <- 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 12
$lastCol = "AY"

# Columns that hold text-like content even though some values look
# numeric or date/time-like (e.g. "Antal" counts stored as text,
# dates/times stored as plain text strings). Force these to Text format
# up front so that assigning the captured values back does not make
# Excel re-interpret them as numbers or dates.
$textCols = @("C","D","F","G","H","I","J","K","L","N","P","T","U","V","W","Y","Z","AA","AB","AC","AT","AW","AX","AY")
foreach ($col in $textCols) {
    $ws.Range("$col$firstRow`:$col$lastRow").NumberFormat = "@"
}

# Snapshot every data row (as it exists before any change) so that
# overwriting one row doesn't corrupt the data needed for another.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:$lastCol$r").Value2
}

# destination row -> source row
$rowMap = @{
    2  = 7
    3  = 11
    4  = 5
    5  = 4
    6  = 3
    7  = 8
    8  = 6
    9  = 10
    10 = 2
    11 = 12
    12 = 9
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("A$destRow`:$lastCol$destRow").Value = $snapshot[$srcRow]
}
